# Apply the re-run analysis R code results: the "count" column (C) for
# "Sheet 1" changed for a handful of country/year rows. Because re-opening
# this particular workbook through the COM layer can re-derive column C
# from its original XML text (which stores single-digit counts as
# space-padded strings, e.g. " 1"), we explicitly re-assert every value in
# C2:C307 via the object model so the saved workbook reflects the correct
# numbers, not just the eleven rows that actually changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

# Row number (1-based) -> corrected "count" value for column C
$counts = [ordered]@{
    2 = 2
    3 = 3
    4 = 4
    5 = 1
    6 = 1
    7 = 1
    8 = 1
    9 = 16
    10 = 19
    11 = 1
    12 = 1
    13 = 7
    14 = 3
    15 = 1
    16 = 6
    17 = 31
    18 = 50
    19 = 13
    20 = 6
    21 = 1
    22 = 1
    23 = 38
    24 = 9
    25 = 1
    26 = 3
    27 = 1
    28 = 1
    29 = 3
    30 = 3
    31 = 1
    32 = 1
    33 = 1
    34 = 3
    35 = 2
    36 = 7
    37 = 36
    38 = 16
    39 = 8
    40 = 16
    41 = 15
    42 = 1
    43 = 12
    44 = 1
    45 = 12
    46 = 2
    47 = 1
    48 = 14
    49 = 2
    50 = 3
    51 = 27
    52 = 6
    53 = 2
    54 = 11
    55 = 15
    56 = 3
    57 = 2
    58 = 1
    59 = 1
    60 = 1
    61 = 2
    62 = 2
    63 = 12
    64 = 4
    65 = 3
    66 = 4
    67 = 7
    68 = 3
    69 = 19
    70 = 27
    71 = 5
    72 = 1
    73 = 2
    74 = 3
    75 = 31
    76 = 1
    77 = 6
    78 = 3
    79 = 1
    80 = 2
    81 = 1
    82 = 3
    83 = 1
    84 = 1
    85 = 2
    86 = 11
    87 = 20
    88 = 1
    89 = 10
    90 = 1
    91 = 2
    92 = 7
    93 = 12
    94 = 1
    95 = 3
    96 = 5
    97 = 1
    98 = 2
    99 = 7
    100 = 3
    101 = 3
    102 = 14
    103 = 3
    104 = 1
    105 = 32
    106 = 1
    107 = 7
    108 = 4
    109 = 1
    110 = 3
    111 = 6
    112 = 1
    113 = 5
    114 = 2
    115 = 5
    116 = 1
    117 = 2
    118 = 1
    119 = 3
    120 = 1
    121 = 5
    122 = 22
    123 = 1
    124 = 1
    125 = 1
    126 = 9
    127 = 7
    128 = 5
    129 = 1
    130 = 2
    131 = 1
    132 = 14
    133 = 1
    134 = 60
    135 = 10
    136 = 2
    137 = 3
    138 = 4
    139 = 57
    140 = 1
    141 = 3
    142 = 1
    143 = 1
    144 = 3
    145 = 13
    146 = 1
    147 = 1
    148 = 1
    149 = 1
    150 = 9
    151 = 3
    152 = 8
    153 = 1
    154 = 2
    155 = 1
    156 = 35
    157 = 2
    158 = 19
    159 = 11
    160 = 16
    161 = 1
    162 = 8
    163 = 1
    164 = 1
    165 = 6
    166 = 1
    167 = 4
    168 = 22
    169 = 1
    170 = 4
    171 = 2
    172 = 1
    173 = 1
    174 = 25
    175 = 2
    176 = 1
    177 = 2
    178 = 2
    179 = 74
    180 = 1
    181 = 6
    182 = 7
    183 = 3
    184 = 6
    185 = 19
    186 = 2
    187 = 1
    188 = 4
    189 = 11
    190 = 1
    191 = 10
    192 = 2
    193 = 1
    194 = 1
    195 = 12
    196 = 10
    197 = 3
    198 = 2
    199 = 1
    200 = 1
    201 = 2
    202 = 3
    203 = 2
    204 = 10
    205 = 1
    206 = 43
    207 = 8
    208 = 2
    209 = 8
    210 = 1
    211 = 3
    212 = 1
    213 = 33
    214 = 4
    215 = 1
    216 = 3
    217 = 1
    218 = 1
    219 = 7
    220 = 9
    221 = 8
    222 = 11
    223 = 2
    224 = 4
    225 = 36
    226 = 17
    227 = 4
    228 = 1
    229 = 9
    230 = 1
    231 = 8
    232 = 1
    233 = 1
    234 = 20
    235 = 2
    236 = 12
    237 = 1
    238 = 2
    239 = 1
    240 = 2
    241 = 8
    242 = 7
    243 = 1
    244 = 1
    245 = 31
    246 = 1
    247 = 1
    248 = 20
    249 = 4
    250 = 1
    251 = 1
    252 = 5
    253 = 2
    254 = 7
    255 = 1
    256 = 7
    257 = 3
    258 = 1
    259 = 1
    260 = 1
    261 = 2
    262 = 10
    263 = 1
    264 = 1
    265 = 1
    266 = 2
    267 = 1
    268 = 2
    269 = 3
    270 = 11
    271 = 1
    272 = 33
    273 = 9
    274 = 1
    275 = 18
    276 = 1
    277 = 2
    278 = 3
    279 = 3
    280 = 1
    281 = 2
    282 = 1
    283 = 1
    284 = 1
    285 = 1
    286 = 5
    287 = 2
    288 = 37
    289 = 14
    290 = 7
    291 = 5
    292 = 1
    293 = 4
    294 = 3
    295 = 7
    296 = 4
    297 = 1
    298 = 1
    299 = 8
    300 = 1
    301 = 1
    302 = 30
    303 = 2
    304 = 7
    305 = 5
    306 = 2
    307 = 9
}

foreach ($row in $counts.Keys) {
    $ws.Cells.Item($row, 3).Value = $counts[$row]
}
